$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived statistics for the Fgf1-Fgfr2 ligand/receptor table,
# and renamed the "Inflammatory-Mac" target cluster to "Resolving-Mac"
# (its row also swapped position with the "MuSCs" row for each sender).

# Row 2
$ws.Range("G2").Value = 0.8775636666666666
$ws.Range("H2").Value = 2.632691
$ws.Range("I2").Value = 0.1887436506618166
$ws.Range("J2").Value = 0.2083714858314108
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2858606666666667
$ws.Range("N2").Value = 0.857582
$ws.Range("O2").Value = 0.0687156860066334
$ws.Range("P2").Value = 0.06932858672617494
$ws.Range("Q2").Value = 0.2508609347957778
$ws.Range("R2").Value = 2.257748413162
$ws.Range("S2").Value = 0.01296964943462309
$ws.Range("T2").Value = 0.01444610062672489

# Row 3
$ws.Range("G3").Value = 0.8775636666666666
$ws.Range("H3").Value = 2.632691
$ws.Range("I3").Value = 0.1887436506618166
$ws.Range("J3").Value = 0.2083714858314108
$ws.Range("O3").Value = 0.90464312565499
$ws.Range("P3").Value = 0.9127119736118995
$ws.Range("Q3").Value = 3.302588293107889
$ws.Range("R3").Value = 29.723294637971
$ws.Range("S3").Value = 0.1707456460822393
$ws.Range("T3").Value = 0.1901831500776309

# Row 4
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 0.8775636666666666
$ws.Range("H4").Value = 2.632691
$ws.Range("I4").Value = 0.1887436506618166
$ws.Range("J4").Value = 0.2083714858314108
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110331
$ws.Range("N4").Value = 0.220662
$ws.Range("O4").Value = 0.02652155835639462
$ws.Range("P4").Value = 0.01783874265571248
$ws.Range("Q4").Value = 0.096822476907
$ws.Range("R4").Value = 0.580934861442
$ws.Range("S4").Value = 0.005005775745426328
$ws.Range("T4").Value = 0.003717085312534975

# Row 5
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.8775636666666666
$ws.Range("H5").Value = 2.632691
$ws.Range("I5").Value = 0.1887436506618166
$ws.Range("J5").Value = 0.2083714858314108
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.0004976666666666667
$ws.Range("N5").Value = 0.001493
$ws.Range("O5").Value = 0.0001196299819817856
$ws.Range("P5").Value = 0.0001206970062130259
$ws.Range("Q5").Value = 0.0004367341847777778
$ws.Range("R5").Value = 0.003930607663
$ws.Range("S5").Value = [double]"2.257939952784955e-05"
$ws.Range("T5").Value = [double]"2.514981452001123e-05"

# Row 6
$ws.Range("G6").Value = 2.458038666666667
$ws.Range("H6").Value = 7.374116000000001
$ws.Range("I6").Value = 0.5286672739959656
$ws.Range("J6").Value = 0.5836444564186148
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2858606666666667
$ws.Range("N6").Value = 0.857582
$ws.Range("O6").Value = 0.0687156860066334
$ws.Range("P6").Value = 0.06932858672617494
$ws.Range("Q6").Value = 0.7026565719457778
$ws.Range("R6").Value = 6.323909147512
$ws.Range("S6").Value = 0.0363277344018896
$ws.Range("T6").Value = 0.04046324531406917

# Row 7
$ws.Range("G7").Value = 2.458038666666667
$ws.Range("H7").Value = 7.374116000000001
$ws.Range("I7").Value = 0.5286672739959656
$ws.Range("J7").Value = 0.5836444564186148
$ws.Range("O7").Value = 0.90464312565499
$ws.Range("P7").Value = 0.9127119736118995
$ws.Range("Q7").Value = 9.25048521593289
$ws.Range("R7").Value = 83.25436694339601
$ws.Range("S7").Value = 0.4782552151792133
$ws.Range("T7").Value = 0.5326992837054783

# Row 8
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 2.458038666666667
$ws.Range("H8").Value = 7.374116000000001
$ws.Range("I8").Value = 0.5286672739959656
$ws.Range("J8").Value = 0.5836444564186148
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.110331
$ws.Range("N8").Value = 0.220662
$ws.Range("O8").Value = 0.02652155835639462
$ws.Range("P8").Value = 0.01783874265571248
$ws.Range("Q8").Value = 0.271197864132
$ws.Range("R8").Value = 1.627187184792
$ws.Range("S8").Value = 0.01402107995840006
$ws.Range("T8").Value = 0.01041148326048487

# Row 9
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 2.458038666666667
$ws.Range("H9").Value = 7.374116000000001
$ws.Range("I9").Value = 0.5286672739959656
$ws.Range("J9").Value = 0.5836444564186148
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.0004976666666666667
$ws.Range("N9").Value = 0.001493
$ws.Range("O9").Value = 0.0001196299819817856
$ws.Range("P9").Value = 0.0001206970062130259
$ws.Range("Q9").Value = 0.001223283909777778
$ws.Range("R9").Value = 0.011009555188
$ws.Range("S9").Value = [double]"6.324445646249708e-05"
$ws.Range("T9").Value = [double]"7.044413858255569e-05"

# Row 10
$ws.Range("G10").Value = 1.313898
$ws.Range("H10").Value = 2.627796
$ws.Range("I10").Value = 0.2825890753422177
$ws.Range("J10").Value = 0.2079840577499744
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2858606666666667
$ws.Range("N10").Value = 0.857582
$ws.Range("O10").Value = 0.0687156860066334
$ws.Range("P10").Value = 0.06932858672617494
$ws.Range("Q10").Value = 0.375591758212
$ws.Range("R10").Value = 2.253550549272
$ws.Range("S10").Value = 0.0194183021701207
$ws.Range("T10").Value = 0.01441924078538088

# Row 11
$ws.Range("G11").Value = 1.313898
$ws.Range("H11").Value = 2.627796
$ws.Range("I11").Value = 0.2825890753422177
$ws.Range("J11").Value = 0.2079840577499744
$ws.Range("O11").Value = 0.90464312565499
$ws.Range("P11").Value = 0.9127119736118995
$ws.Range("Q11").Value = 4.944671615246
$ws.Range("R11").Value = 29.668029691476
$ws.Range("S11").Value = 0.2556422643935373
$ws.Range("T11").Value = 0.1898295398287904

# Row 12
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = 1.313898
$ws.Range("H12").Value = 2.627796
$ws.Range("I12").Value = 0.2825890753422177
$ws.Range("J12").Value = 0.2079840577499744
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.110331
$ws.Range("N12").Value = 0.220662
$ws.Range("O12").Value = 0.02652155835639462
$ws.Range("P12").Value = 0.01783874265571248
$ws.Range("Q12").Value = 0.144963680238
$ws.Range("R12").Value = 0.579854720952
$ws.Range("S12").Value = 0.007494702652568222
$ws.Range("T12").Value = 0.003710174082692635

# Row 13
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("G13").Value = 1.313898
$ws.Range("H13").Value = 2.627796
$ws.Range("I13").Value = 0.2825890753422177
$ws.Range("J13").Value = 0.2079840577499744
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.0004976666666666667
$ws.Range("N13").Value = 0.001493
$ws.Range("O13").Value = 0.0001196299819817856
$ws.Range("P13").Value = 0.0001206970062130259
$ws.Range("Q13").Value = 0.000653883238
$ws.Range("R13").Value = 0.003923299428
$ws.Range("S13").Value = [double]"3.380612599143896e-05"
$ws.Range("T13").Value = [double]"2.5103053110459e-05"
